$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Edit row 132: re-coded by "chen" ---
$ws.Cells.Item(132, 6).Value = "3: 4496"
$ws.Cells.Item(132, 9).Value = " trimethoprim"
$ws.Cells.Item(132, 11).Value = 0.04008283786492084
$ws.Cells.Item(132, 12).Value = "chen"
$ws.Cells.Item(132, 13).Value = "1/31/19 13:50:08"

# --- Append 12 new coded segments (rows 155-166) ---
$newRows = @(
    @{ E="Drug Resisted"; F="3: 4362"; G="3: 4371"; I="cefotaxime";                    J=10; K=0.03340236488743403;  M="1/31/19 13:46:42" },
    @{ E="Drug Resisted"; F="3: 4325"; G="3: 4337"; I="ciprofloxacin";                 J=13; K=0.043423074353664236; M="1/31/19 13:46:48" },
    @{ E="Drug Resisted"; F="3: 4303"; G="3: 4317"; I="chloramphenicol";               J=15; K=0.05010354733115105;  M="1/31/19 13:47:04" },
    @{ E="Drug Resisted"; F="3: 4497"; G="3: 4508"; I="trimethoprim";                  J=12; K=0.04008283786492084;  M="1/31/19 13:51:58" },
    @{ E="Drug Resisted"; F="3: 4478"; G="3: 4489"; I="tetracycline";                  J=12; K=0.04008283786492084;  M="1/31/19 13:47:31" },
    @{ E="Drug Resisted"; F="3: 4287"; G="3: 4295"; I="ceftiofur";                     J=9;  K=0.030062128398690626; M="1/31/19 13:47:40" },
    @{ E="Drug Resisted"; F="3: 4270"; G="3: 4279"; I="ampicillin";                    J=10; K=0.03340236488743403;  M="1/31/19 13:47:44" },
    @{ E="Drug Resisted"; F="3: 4236"; G="3: 4262"; I="amoxicillin/clavulanic acid";   J=27; K=0.09018638519607187;  M="1/31/19 13:47:52" },
    @{ E="Drug Resisted"; F="3: 4416"; G="3: 4428"; I="spectinomycin";                 J=13; K=0.043423074353664236; M="1/31/19 13:47:58" },
    @{ E="Drug Resisted"; F="3: 4400"; G="3: 4407"; I="neomycin";                      J=8;  K=0.026721891909947226; M="1/31/19 13:48:14" },
    @{ E="Drug Resisted"; F="3: 4345"; G="3: 4354"; I="gentamicin";                    J=10; K=0.03340236488743403;  M="1/31/19 13:48:18" },
    @{ E="Drug Resisted"; F="3: 4436"; G="3: 4447"; I="streptomycin";                  J=12; K=0.04008283786492084;  M="1/31/19 13:48:49" }
)

$r = 155
foreach ($row in $newRows) {
    $ws.Range("A132:M132").Copy()
    $ws.Range("A$r`:M$r").PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = "●"
    $ws.Cells.Item($r, 4).Value = "20373"
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = 0
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = "chen"
    $ws.Cells.Item($r, 13).Value = $row.M

    $r = $r + 1
}
